$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell H12: a short "essence of the market" headline next to the title in D12.
$cell = $ws.Range("H12")
$cell.Value = "Lĩnh ngộ tinh túy thị trường"

# Big bold emphasis font for the new headline (Calibri 20 bold, matches the
# sheet's default theme font family/scheme).
$cell.Style = "Normal"
$cell.Font.Bold = $true
$cell.Font.Size = 20

# Row 12 grows taller to fit the larger font.
$ws.Rows(12).RowHeight = 26.25

# Move the active selection to the newly-added cell.
$cell.Select()
